# Natmi following Dr Hou advice
# Rewrites the LR-pair result rows for Fgf10-Fgfrl1 (sheet1) with the
# updated per-target-cluster statistics, and adds two new target-cluster
# rows (ECs -> M2 / M2 already covered, sCs) that weren't present before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster FAPs / Fgf10 / Fgfrl1 -> Target cluster ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf10"
$ws.Range("C2").Value = "Fgfrl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.422753333333333
$ws.Range("H2").Value = 4.26826
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.6019246666666667
$ws.Range("N2").Value = 1.805774
$ws.Range("O2").Value = 0.05746655956902961
$ws.Range("P2").Value = 0.0574665595690296
$ws.Range("Q2").Value = 0.8563903259155555
$ws.Range("R2").Value = 7.707512933239999
$ws.Range("S2").Value = 0.05746655956902961
$ws.Range("T2").Value = 0.0574665595690296

# Row 3: Sending cluster FAPs / Fgf10 / Fgfrl1 -> Target cluster FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf10"
$ws.Range("C3").Value = "Fgfrl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.422753333333333
$ws.Range("H3").Value = 4.26826
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.794889
$ws.Range("N3").Value = 17.384667
$ws.Range("O3").Value = 0.5532458667270895
$ws.Range("P3").Value = 0.5532458667270894
$ws.Range("Q3").Value = 8.244697641046667
$ws.Range("R3").Value = 74.20227876941999
$ws.Range("S3").Value = 0.5532458667270895
$ws.Range("T3").Value = 0.5532458667270894

# Row 4 (new): Sending cluster FAPs / Fgf10 / Fgfrl1 -> Target cluster M2
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf10"
$ws.Range("C4").Value = "Fgfrl1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.422753333333333
$ws.Range("H4").Value = 4.26826
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3673473333333333
$ws.Range("N4").Value = 1.102042
$ws.Range("O4").Value = 0.03507114524883653
$ws.Range("P4").Value = 0.03507114524883652
$ws.Range("Q4").Value = 0.522644642991111
$ws.Range("R4").Value = 4.70380178692
$ws.Range("S4").Value = 0.03507114524883653
$ws.Range("T4").Value = 0.03507114524883652

# Row 5 (new): Sending cluster FAPs / Fgf10 / Fgfrl1 -> Target cluster sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf10"
$ws.Range("C5").Value = "Fgfrl1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.422753333333333
$ws.Range("H5").Value = 4.26826
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.710185666666666
$ws.Range("N5").Value = 11.130557
$ws.Range("O5").Value = 0.3542164284550445
$ws.Range("P5").Value = 0.3542164284550444
$ws.Range("Q5").Value = 5.278679024535555
$ws.Range("R5").Value = 47.50811122082
$ws.Range("S5").Value = 0.3542164284550445
$ws.Range("T5").Value = 0.3542164284550444
